$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from serial date 46073 (2026-02-20) to 46074 (2026-02-21).
for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46073) {
        $cell.Value = 46074
    }
}
